$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A61").Value = "2025/12/05 06:00"
$ws.Range("B61").Value = "-"
$ws.Range("C61").Value = "-"
$ws.Range("D61").Value = "-"
$ws.Range("E61").Value = "-"
$ws.Range("F61").Value = "-"
$ws.Range("G61").Value = "-"
